$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "edit2"
$ws.Range("B5").Value = "riya-morankar"
$ws.Range("C5").Value = "Merged"
# (D5 intentionally left blank - no "Comment" for this log entry)

# E5's value looks like a date ("2025-06-18"), so Excel would normally
# auto-convert it into a date serial number. Format the cell as text
# first so the literal string is preserved, then clear the formatting
# again afterwards so the cell doesn't end up with a lingering explicit
# style (matching the plain, unstyled cells used elsewhere in this sheet).
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2025-06-18"
$ws.Range("E5").ClearFormats()

$ws.Range("F5").Value = "N/A"
